$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "59.774.56"
Set-TextValue "E2" "  -0.58%  "
Set-TextValue "D3" "2.649.08"
Set-TextValue "E3" "  +1.61%  "
Set-TextValue "E4" "  +0.14%  "
Set-TextValue "D5" "517.74"
Set-TextValue "E5" "  -0.36%  "
Set-TextValue "D6" "146.61"
Set-TextValue "E6" "  -1.35%  "
Set-TextValue "D7" "0.996"
Set-TextValue "E7" "  -0.35%  "
Set-TextValue "E8" "  +0.59%  "
Set-TextValue "D9" "2.677.99"
Set-TextValue "E9" "  +2.44%  "
Set-TextValue "D10" "6.49"
Set-TextValue "E10" "  +3.32%  "
Set-TextValue "E11" "  +2.16%  "
Set-TextValue "D12" "0.340"
Set-TextValue "E12" "  -0.22%  "
Set-TextValue "E13" "  -1.52%  "
Set-TextValue "D14" "3.118.08"
Set-TextValue "E14" "  +1.75%  "
Set-TextValue "D15" "59.740.64"
Set-TextValue "E15" "  -0.74%  "
Set-TextValue "D16" "21.31"
Set-TextValue "E16" "  +0.68%  "
Set-TextValue "E17" "  +1.07%  "
Set-TextValue "D18" "2.692.27"
Set-TextValue "E18" "  +3.50%  "
Set-TextValue "E19" "  +0.23%  "
Set-TextValue "D20" "347.06"
Set-TextValue "E20" "  +1.49%  "
Set-TextValue "D21" "10.57"
Set-TextValue "E21" "  +1.81%  "
Set-TextValue "E22" "  +1.74%  "
Set-TextValue "D23" "0.998"
Set-TextValue "E23" "  +0.14%  "
Set-TextValue "D24" "61.20"
Set-TextValue "E24" "  +1.59%  "
Set-TextValue "D25" "0.424"
Set-TextValue "E25" "  +1.50%  "
Set-TextValue "D26" "2.768.54"
Set-TextValue "E26" "  +1.54%  "
Set-TextValue "B27" "Kaspa"
Set-TextValue "C27" "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D27" "0.162"
Set-TextValue "E27" "  +0.58%  "
Set-TextValue "B28" "Binance-PegBSC-USD"
Set-TextValue "C28" "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-TextValue "D28" "0.992"
Set-TextValue "E28" "  -0.67%  "
Set-TextValue "D29" "0.0₃0823"
Set-TextValue "E29" "  +2.41%  "
Set-TextValue "D30" "7.22"
Set-TextValue "E30" "  +2.35%  "
Set-TextValue "D31" "0.998"
Set-TextValue "E31" "  -0.23%  "
Set-TextValue "D32" "6.52"
Set-TextValue "E32" "  +8.61%  "
Set-TextValue "D33" "19.09"
Set-TextValue "E33" "  +0.94%  "
Set-TextValue "E34" "  +0.20%  "
Set-TextValue "D35" "1.07"
Set-TextValue "E35" "  +17.89%  "
Set-TextValue "D36" "149.84"
Set-TextValue "E36" "  +0.03%  "
Set-TextValue "D37" "4.07"
Set-TextValue "E37" "  +3.02%  "
Set-TextValue "D38" "1.17"
Set-TextValue "E38" "  +2.81%  "
Set-TextValue "D39" "0.875"
Set-TextValue "E39" "  +1.69%  "
Set-TextValue "D40" "36.73"
Set-TextValue "E40" "  +0.90%  "
Set-TextValue "E41" "  +3.54%  "
Set-TextValue "E42" "  +0.72%  "
Set-TextValue "D43" "283.31"
Set-TextValue "E43" "  -1.27%  "
Set-TextValue "D44" "0.623"
Set-TextValue "E44" "  -0.05%  "
Set-TextValue "D45" "0.0994"
Set-TextValue "E45" "  -0.88%  "
Set-TextValue "B46" "FirstDigitalUSD"
Set-TextValue "C46" "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue "D46" "0.994"
Set-TextValue "E46" "  -0.55%  "
Set-TextValue "B47" "EnergySwap"
Set-TextValue "C47" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D47" "19.88"
Set-TextValue "E47" "  +2.22%  "
Set-TextValue "E48" "  +0.06%  "
Set-TextValue "E49" "  +1.47%  "
Set-TextValue "E50" "  +2.45%  "
Set-TextValue "B51" "Maker"
Set-TextValue "C51" "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue "D51" "2.000.55"
Set-TextValue "E51" "  +2.48%  "
